# Update the "timestamp" column (Z) of the Log_Muestras sheet.
# The logging run was re-executed, producing a new set of timestamps for
# each previously-logged sample row. Rows sharing the same timestamp in
# the original run are updated in contiguous blocks to the corresponding
# new timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2:Z6").Value   = "2025-10-17T07:09:34.280053"
$ws.Range("Z7:Z45").Value  = "2025-10-17T07:09:34.280562"
$ws.Range("Z46:Z64").Value = "2025-10-17T07:09:34.330215"
$ws.Range("Z65:Z66").Value = "2025-10-17T07:09:34.341596"
$ws.Range("Z67:Z74").Value = "2025-10-17T07:09:34.342133"
$ws.Range("Z75:Z102").Value = "2025-10-17T07:09:34.392257"
$ws.Range("Z103:Z105").Value = "2025-10-17T07:09:34.461058"
$ws.Range("Z106:Z109").Value = "2025-10-17T07:09:34.462059"
$ws.Range("Z110:Z112").Value = "2025-10-17T07:09:34.463059"
